# Fix misspelled "SD\Foilo.xlsx" -> "SD\Folio.xlsx" wherever it appears in the
# TestDataFileName column of the TestDataMappingSheet_SD sheet, and restore
# the last scroll position / selection that Excel saved on exit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestDataMappingSheet_SD")

# Replace every occurrence of the misspelled file name with the corrected one.
$ws.Cells.Replace("SD\Foilo.xlsx", "SD\Folio.xlsx")

# Restore the saved view state (scroll position + active cell selection).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B45").Select()
